$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=4; Value=-12.04570000000001},
    @{Row=6; Value=-12.35360000000001},
    @{Row=7; Value=-13.1336},
    @{Row=16; Value=-14.26439999999998},
    @{Row=20; Value=-11.5477},
    @{Row=28; Value=-12.70779999999999},
    @{Row=29; Value=-11.24030000000001},
    @{Row=32; Value=-12.9231},
    @{Row=40; Value=-12.7645},
    @{Row=46; Value=-14.58299999999999},
    @{Row=51; Value=-11.5427},
    @{Row=52; Value=-11.2672},
    @{Row=57; Value=-14.16169999999999},
    @{Row=59; Value=-12.61389999999999},
    @{Row=62; Value=-14.1528},
    @{Row=66; Value=-10.9659},
    @{Row=73; Value=-12.3742},
    @{Row=74; Value=-11.79380000000001},
    @{Row=92; Value=-10.92679999999999},
    @{Row=100; Value=-13.0562}
)

foreach ($u in $updates) {
    $ws.Range("C$($u.Row)").Value = $u.Value
}
